$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G10").Value = "arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G18").Value = "aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G19").Value = "Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G24").Value = "shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G32").Value = "arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G40").Value = "aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G41").Value = "Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G46").Value = "hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("G52").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G54").Value = "yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G58").Value = "afaf.abdallah@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg"
$ws.Range("G62").Value = "wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G68").Value = "hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("G74").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G76").Value = "yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G80").Value = "afaf.abdallah@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg"
$ws.Range("G84").Value = "wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G90").Value = "shaimaa.ahmed@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G96").Value = "aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G98").Value = "yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G106").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G107").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G112").Value = "shaimaa.ahmed@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G118").Value = "aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G120").Value = "yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G128").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G129").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G134").Value = "majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G142").Value = "yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, merna.said@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"
$ws.Range("G150").Value = "Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G156").Value = "majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"
$ws.Range("G164").Value = "yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, merna.said@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"
$ws.Range("G172").Value = "Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
